# Generate Report for Archive
# - Status text moves from "Ready for handoff" to "In Translation" for every
#   file row (Overview sheet columns E/F, and the per-locale "Status" column
#   C on the zh-cn / de-de sheets).
# - Because the status text got shorter, the report's own column-sizing
#   logic narrows the two/three affected "Status" columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn status) and F (de-de status) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

# --- zh-cn sheet: column C ("Status") ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

# --- de-de sheet: column C ("Status") ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"

# --- Re-size the Status columns now that the text is shorter ---
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C
